$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 30
$ws.Cells.Item(14, 3).Value = 0.001
$ws.Cells.Item(14, 4).Value = 0.003
$ws.Cells.Item(14, 5).Value = 'Regular'
$ws.Cells.Item(14, 6).Value = '<function relu at 0x11ad159d8>'
$ws.Cells.Item(14, 7).Value = 0.9657999873161316
$ws.Cells.Item(14, 8).Value = 0.2011000066995621
$ws.Cells.Item(14, 9).Value = 0.1811999976634979
$ws.Cells.Item(14, 10).Value = 0.1451183259487152
$ws.Cells.Item(14, 11).Value = 5.327450752258301
$ws.Cells.Item(14, 12).Value = 0.2011000066995621
$ws.Cells.Item(14, 13).Value = 'logs/results_71.log'
$ws.Cells.Item(14, 14).Value = 'weights/model_71.ckpt'
$ws.Cells.Item(14, 15).Value = 'tb/71/non_robust'
$ws.Cells.Item(14, 16).Value = '(5.475276, 12.719564, 18.869154, 27.198263, 26.215324, 22.13533, 16.86695)'
$ws.Cells.Item(14, 17).Value = '(78.487724, 11.971958, 8.304008, 4.5321946, 2.381914, 1.7097418, 1.7520251, 2.2926486)'

# Row 15
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 30
$ws.Cells.Item(15, 3).Value = 0.005
$ws.Cells.Item(15, 4).Value = 0.003
$ws.Cells.Item(15, 5).Value = 'Regular'
$ws.Cells.Item(15, 6).Value = '<function relu at 0x12008f9d8>'
$ws.Cells.Item(15, 7).Value = 0.9585999846458435
$ws.Cells.Item(15, 8).Value = 0.2369000017642975
$ws.Cells.Item(15, 9).Value = 0.03700000047683716
$ws.Cells.Item(15, 10).Value = 0.1429557055234909
$ws.Cells.Item(15, 11).Value = 4.122503280639648
$ws.Cells.Item(15, 12).Value = 0.2369000017642975
$ws.Cells.Item(15, 13).Value = 'logs/results_72.log'
$ws.Cells.Item(15, 14).Value = 'weights/model_72.ckpt'
$ws.Cells.Item(15, 15).Value = 'tb/72/non_robust'
$ws.Cells.Item(15, 16).Value = '(2.139943, 3.496155, 3.843952, 5.1138744, 6.2203045, 6.9547405, 8.452685)'
$ws.Cells.Item(15, 17).Value = '(29.190838, 5.508862, 3.9576876, 1.9097207, 1.488211, 1.3313731, 1.552425, 1.4604229)'

# Row 16
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 70
$ws.Cells.Item(16, 3).Value = 0.008
$ws.Cells.Item(16, 4).Value = 0.003
$ws.Cells.Item(16, 5).Value = 'Regular'
$ws.Cells.Item(16, 6).Value = '<function relu at 0x1240199d8>'
$ws.Cells.Item(16, 7).Value = 0.9384999871253967
$ws.Cells.Item(16, 8).Value = 0.2345000058412552
$ws.Cells.Item(16, 9).Value = 0.07829999923706055
$ws.Cells.Item(16, 10).Value = 0.2154168486595154
$ws.Cells.Item(16, 11).Value = 4.516582012176514
$ws.Cells.Item(16, 12).Value = 0.2345000058412552
$ws.Cells.Item(16, 13).Value = 'logs/results_74.log'
$ws.Cells.Item(16, 14).Value = 'weights/model_74.ckpt'
$ws.Cells.Item(16, 15).Value = 'tb/74/non_robust'
$ws.Cells.Item(16, 16).Value = '(1.7076006, 2.5921235, 2.5732467, 3.4120953, 4.31576, 5.6463065, 6.7571893)'
$ws.Cells.Item(16, 17).Value = '(23.436884, 4.6997, 2.8007762, 1.6223694, 1.4685857, 1.5694603, 1.4128212, 1.4283097)'

# Row 17
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 30
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 0.003
$ws.Cells.Item(17, 5).Value = 'FGSM'
$ws.Cells.Item(17, 6).Value = '<function relu at 0x11951f9d8>'
$ws.Cells.Item(17, 7).Value = 0.9517999887466431
$ws.Cells.Item(17, 8).Value = 0.8098999857902527
$ws.Cells.Item(17, 9).Value = 0.4519999921321869
$ws.Cells.Item(17, 10).Value = 0.1830078810453415
$ws.Cells.Item(17, 11).Value = 0.6702156066894531
$ws.Cells.Item(17, 12).Value = 0.8098999857902527
$ws.Cells.Item(17, 13).Value = 'logs/results_75.log'
$ws.Cells.Item(17, 14).Value = 'weights/model_75.ckpt'
$ws.Cells.Item(17, 15).Value = 'tb/75/robust'
$ws.Cells.Item(17, 16).Value = '(6.767173, 17.473831, 38.257633, 49.1879, 27.285254, 10.4143715, 4.247443)'
$ws.Cells.Item(17, 17).Value = '(131.64798, 17.081043, 17.57381, 15.574147, 13.246862, 13.433027, 11.568809, 11.362484)'
